$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (min) - inline string cleanup / addition
$ws.Range("AZ3").Value = ""
$ws.Range("BI3").Value = ""
$ws.Range("BR3").Value = ""
$ws.Range("CB3").Value = "-inf"
$ws.Range("CZ3").Value = ""

# Row 5
$ws.Range("L5").Value = 0.0002499296011591896
$ws.Range("AP5").Value = -0.4289874419420469
$ws.Range("AQ5").Value = -0.002515802716733498
$ws.Range("AR5").Value = -0.1992363246562166
$ws.Range("AS5").Value = -0.3882963372003823
$ws.Range("AT5").Value = -0.468523624718505
$ws.Range("AU5").Value = -0.4533567389337844
$ws.Range("AV5").Value = -0.4213262799231051
$ws.Range("AW5").Value = -0.4677501893286136
$ws.Range("AX5").Value = -0.392145303242806
$ws.Range("AY5").Value = -0.5043743037254131
$ws.Range("AZ5").Value = -0.4023202441444778
$ws.Range("BA5").Value = -0.4540274990402154
$ws.Range("BB5").Value = -0.4087074000922734
$ws.Range("BC5").Value = -0.4763942240812771
$ws.Range("BD5").Value = -0.4272851641602656
$ws.Range("BE5").Value = -0.4313220124072843
$ws.Range("BF5").Value = -0.458103386079036
$ws.Range("BG5").Value = -0.4150513950716801
$ws.Range("BH5").Value = -0.5058346683304387
$ws.Range("BI5").Value = -0.4085989974825165
$ws.Range("BJ5").Value = -0.4584826887915781
$ws.Range("BK5").Value = -0.4034574354015902
$ws.Range("BL5").Value = -0.4811782871077419
$ws.Range("BM5").Value = -0.4330943822783881
$ws.Range("BN5").Value = -0.4238386894340955
$ws.Range("BO5").Value = -0.4672460622804273
$ws.Range("BP5").Value = -0.3979751603767134
$ws.Range("BQ5").Value = -0.505228716615331
$ws.Range("BR5").Value = -0.4109879125565294
$ws.Range("BS5").Value = -0.4548635967006624
$ws.Range("BV5").Value = -0.3976285032279087
$ws.Range("BX5").Value = -0.4034585177949079
$ws.Range("BY5").Value = -0.002163661202950253
$ws.Range("BZ5").Value = -0.003687601923746804
$ws.Range("CA5").Value = -0.002175536777061334
$ws.Range("CB5").Value = -0.003801839157656842
$ws.Range("CX5").Value = -0.1099631785882957
$ws.Range("CY5").Value = -0.1179916993691728
$ws.Range("CZ5").Value = -0.1284488991882546
$ws.Range("DA5").Value = -0.1198677298574224
$ws.Range("DB5").Value = -0.1390576034388533
$ws.Range("DC5").Value = -0.1310882107956524

# Row 6
$ws.Range("L6").Value = -0.0007976935587679249
$ws.Range("AP6").Value = -0.006735350774594183
$ws.Range("AQ6").Value = -0.0001273410084810027
$ws.Range("AR6").Value = -0.0004235379357077261
$ws.Range("AS6").Value = -0.1226600649984384
$ws.Range("AT6").Value = -0.09776873637607378
$ws.Range("AU6").Value = 0.0003818244244859251
$ws.Range("AV6").Value = -0.3138266119853111
$ws.Range("AW6").Value = -0.01580312239465077
$ws.Range("AX6").Value = -0.1064646296983257
$ws.Range("AY6").Value = -0.1056516540922313
$ws.Range("AZ6").Value = -0.1534318906195552
$ws.Range("BA6").Value = -0.3048655654032916
$ws.Range("BB6").Value = -0.09409491640815729
$ws.Range("BC6").Value = -0.0905350884807213
$ws.Range("BD6").Value = 0.001496611916891348
$ws.Range("BE6").Value = -0.3024532882813517
$ws.Range("BF6").Value = -0.01941555703679932
$ws.Range("BG6").Value = -0.07324860731844889
$ws.Range("BH6").Value = -0.1026070102020353
$ws.Range("BI6").Value = -0.1524247129058478
$ws.Range("BJ6").Value = -0.3135947887161432
$ws.Range("BK6").Value = -0.1053970773515126
$ws.Range("BL6").Value = -0.08102454940681041
$ws.Range("BM6").Value = 0.0002178050055453847
$ws.Range("BN6").Value = -0.2980252365156191
$ws.Range("BO6").Value = -0.01506823089642617
$ws.Range("BP6").Value = -0.1001911019084227
$ws.Range("BQ6").Value = -0.1027235533736502
$ws.Range("BR6").Value = -0.140680404465216
$ws.Range("BS6").Value = -0.3090546966650163
$ws.Range("BV6").Value = -0.1448653693260238
$ws.Range("BX6").Value = 0.03392884867741811
$ws.Range("BY6").Value = -0.003919275009171727
$ws.Range("BZ6").Value = -0.002636397710239916
$ws.Range("CA6").Value = -0.00392401074007146
$ws.Range("CB6").Value = -0.002454426134611529
$ws.Range("CX6").Value = -0.03257885628080859
$ws.Range("CY6").Value = -0.005843192745412799
$ws.Range("CZ6").Value = -0.004003829119099353
$ws.Range("DA6").Value = -0.001320786793736566
$ws.Range("DB6").Value = 0.00101197522859596
$ws.Range("DC6").Value = 0.001512186210608543

# Row 7
$ws.Range("AP7").Value = -0.02427073914997907
$ws.Range("AR7").Value = -0.03993965114976979
$ws.Range("AS7").Value = 0.004288958099441264
$ws.Range("AT7").Value = -0.04855721954005209
$ws.Range("AU7").Value = -0.0007758462001577733
$ws.Range("AW7").Value = 0
$ws.Range("AX7").Value = 0.04361624429993665
$ws.Range("AY7").Value = 0.1378849372678146
$ws.Range("AZ7").Value = 0.02835378085776679
$ws.Range("BA7").Value = 1.513490551321156
$ws.Range("BB7").Value = 0.04074243736388703
$ws.Range("BC7").Value = -0.01122996754147931
$ws.Range("BD7").Value = -0.0008799200848481125
$ws.Range("BF7").Value = -0.03736427528291004
$ws.Range("BG7").Value = 0.1089242135336597
$ws.Range("BH7").Value = 0.2052041887598985
$ws.Range("BI7").Value = 0.03675936343979807
$ws.Range("BJ7").Value = 8.954698462895109
$ws.Range("BK7").Value = 0.05995903871408178
$ws.Range("BL7").Value = -0.01551853080425364
$ws.Range("BM7").Value = -0.0004869054805961527
$ws.Range("BO7").Value = 0
$ws.Range("BP7").Value = 0.05843233003931114
$ws.Range("BQ7").Value = 0.1582300068652973
$ws.Range("BR7").Value = 0.05942655830956672
$ws.Range("BS7").Value = 1.883871438809595
$ws.Range("BV7").Value = 0.06884279354814198
$ws.Range("DC7").Value = 0.8001109520989129

# Row 8
$ws.Range("L8").Value = 0.0004091116079876026
$ws.Range("AP8").Value = -0.5332016732343132
$ws.Range("AQ8").Value = -0.002958574285957084
$ws.Range("AR8").Value = -0.2088720741693975
$ws.Range("AS8").Value = -0.439181873533821
$ws.Range("AT8").Value = -0.5345929437447041
$ws.Range("AU8").Value = -0.4785844857453802
$ws.Range("AV8").Value = -0.4737453317713382
$ws.Range("AW8").Value = -0.4859435487232779
$ws.Range("AX8").Value = -0.446780930951804
$ws.Range("AY8").Value = -0.5703786491769666
$ws.Range("AZ8").Value = -0.4590009444183961
$ws.Range("BA8").Value = -0.5154139307013563
$ws.Range("BB8").Value = -0.4475283354631436
$ws.Range("BC8").Value = -0.5392697016609356
$ws.Range("BD8").Value = -0.4781088587427645
$ws.Range("BE8").Value = -0.4791246759204425
$ws.Range("BF8").Value = -0.4859428213602781
$ws.Range("BG8").Value = -0.4593684882039413
$ws.Range("BH8").Value = -0.5703384422200263
$ws.Range("BI8").Value = -0.4626841542747325
$ws.Range("BJ8").Value = -0.5166560500015543
$ws.Range("BK8").Value = -0.448277049114476
$ws.Range("BL8").Value = -0.5376821941195509
$ws.Range("BM8").Value = -0.4780984567448335
$ws.Range("BN8").Value = -0.4741086389774381
$ws.Range("BO8").Value = -0.4859545940190099
$ws.Range("BP8").Value = -0.4500477555439378
$ws.Range("BQ8").Value = -0.5707262361356821
$ws.Range("BR8").Value = -0.4629867670887998
$ws.Range("BS8").Value = -0.5153204865687687
$ws.Range("BV8").Value = -0.4469134273465155
$ws.Range("BX8").Value = -0.4611992006570705
$ws.Range("BY8").Value = -0.002644185247627125
$ws.Range("BZ8").Value = -0.003813509228644545
$ws.Range("CA8").Value = -0.002651656495131854
$ws.Range("CB8").Value = -0.003855303775955623
$ws.Range("CX8").Value = -0.1428177754500709
$ws.Range("CY8").Value = -0.1441326273215782
$ws.Range("CZ8").Value = -0.1437556542526778
$ws.Range("DA8").Value = -0.1419865922419628
$ws.Range("DB8").Value = -0.1441244609606021
$ws.Range("DC8").Value = -0.1428756793141632

# Row 9
$ws.Range("L9").Value = 0.0002045348867338558
$ws.Range("AP9").Value = -0.3167735904067475
$ws.Range("AQ9").Value = -0.001480382909758353
$ws.Range("AR9").Value = -0.1105462767346451
$ws.Range("AS9").Value = -0.2511220884108152
$ws.Range("AT9").Value = -0.3177925123136687
$ws.Range("AU9").Value = -0.2779089293900528
$ws.Range("AV9").Value = -0.2745658760241136
$ws.Range("AW9").Value = -0.2830226982137287
$ws.Range("AX9").Value = -0.2562130217271911
$ws.Range("AY9").Value = -0.3445449284481557
$ws.Range("AZ9").Value = -0.2644736200641042
$ws.Range("BA9").Value = -0.3038778345012682
$ws.Range("BB9").Value = -0.2567156233736266
$ws.Range("BC9").Value = -0.3212288321245043
$ws.Range("BD9").Value = -0.2775796644215811
$ws.Range("BE9").Value = -0.2782830720569474
$ws.Range("BF9").Value = -0.2830221909712115
$ws.Range("BG9").Value = -0.2647235133665304
$ws.Range("BH9").Value = -0.3445142581413584
$ws.Range("BI9").Value = -0.2669816880014064
$ws.Range("BJ9").Value = -0.3047705774361634
$ws.Range("BK9").Value = -0.2572194463466858
$ws.Range("BL9").Value = -0.3200604395385945
$ws.Range("BM9").Value = -0.2775724650491466
$ws.Range("BN9").Value = -0.2748163260093607
$ws.Range("BO9").Value = -0.2830304009366993
$ws.Range("BP9").Value = -0.2584123487705164
$ws.Range("BQ9").Value = -0.3448101314395056
$ws.Range("BR9").Value = -0.2671881326621408
$ws.Range("BS9").Value = -0.3038107201118138
$ws.Range("BV9").Value = -0.2563020958389861
$ws.Range("BX9").Value = -0.2659694833708006
$ws.Range("BY9").Value = -0.001322967745641519
$ws.Range("BZ9").Value = -0.001908575945391942
$ws.Range("CA9").Value = -0.001326708325055777
$ws.Range("CB9").Value = -0.001929513398956261
$ws.Range("CX9").Value = -0.07415863964179641
$ws.Range("CY9").Value = -0.07486899701803211
$ws.Range("CZ9").Value = -0.07466527907609444
$ws.Range("DA9").Value = -0.07370986847638444
$ws.Range("DB9").Value = -0.07486458340446293
$ws.Range("DC9").Value = -0.07418991111252371
